$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.481.88"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.07%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.849.90"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.58%  "
# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.10%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.92"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.24%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6264"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.18%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.03%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.31"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.64%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07537"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.12%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.2976"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.17%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "24.31"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.42%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.977.45"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +6.31%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07688"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.18%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.003"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.90%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6854"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.11%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.74"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.34%  "
# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.51%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.199.15"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +4.20%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.222"
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "29.615.58"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.66%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "234.50"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.93%  "
# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.41%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9995"
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.594"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.16%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.000"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.08%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.35"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.27%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1391"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.13%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.433"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.57%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "17.72"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.26%  "
# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.24%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05870"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -5.41%  "
# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.65%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.099"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.52%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.036"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.04%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.886"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.77%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.171"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.44%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.7190"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.64%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.588"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.80%  "
# Row 39
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.238.02"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.06%  "
# Row 40
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.791"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.00%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01777"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.58%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9094"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.49%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.135"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.52%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.106.18"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +4.23%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9996"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.04%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.95"
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "67.25"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.66%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.320"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +8.90%  "
# Row 49
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.714"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.86%  "
# Row 50
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4029"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.79%  "
# Row 51
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.132"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.05%  "
